# Update generator linear examples: replace the numeric (text) values that
# were regenerated for this run across the follower-restrictions, modified
# point and bf/BF vector sheets.
#
# The workbook stores these numbers as plain text (shared strings) rather
# than numeric cells, so each new value is entered with a leading
# apostrophe (forces text entry) and the cell style is reset back to
# "Normal" right afterwards so no stray quote-prefix number format lingers
# on the cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Worksheet,
        [string]$Cell,
        [string]$Value
    )
    $range = $Worksheet.Range($Cell)
    $range.Value = "'" + $Value
    $range.Style = "Normal"
}

# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively and would hit the same sheet for both, so every sheet
# below is looked up by its 1-based tab position instead:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_follower ---
$wsFollower = $wb.Worksheets.Item(3)

Set-TextValue $wsFollower "B2" "-4.374623078112156"
Set-TextValue $wsFollower "D2" "0.33468162538227564"
Set-TextValue $wsFollower "E2" "0.25385517675070546"
Set-TextValue $wsFollower "F2" "0.6571555200523318"

Set-TextValue $wsFollower "B3" "-2.8102693382873367"
Set-TextValue $wsFollower "D3" "0.9092567913461869"
Set-TextValue $wsFollower "E3" "0.6522943366696484"
Set-TextValue $wsFollower "F3" "0.7849656311840086"

Set-TextValue $wsFollower "B4" "0.9341385726238034"
Set-TextValue $wsFollower "D4" "0.7906785535517057"
Set-TextValue $wsFollower "E4" "0.5648689493855065"
Set-TextValue $wsFollower "F4" "0.13732970979821657"

Set-TextValue $wsFollower "B5" "0.36494658748581443"
Set-TextValue $wsFollower "D5" "0.5618257705012442"
Set-TextValue $wsFollower "E5" "0.0021050482917261888"
Set-TextValue $wsFollower "F5" "0.8770541189264485"

# --- Punto_modificado ---
$wsPunto = $wb.Worksheets.Item(4)

Set-TextValue $wsPunto "A2" "4.184892416399492"
Set-TextValue $wsPunto "B2" "4.374623078112156"

# --- Vector_bf ---
$wsBf = $wb.Worksheets.Item(5)

Set-TextValue $wsBf "A2" "-2.5941065025660786"

# --- Vector_BF ---
$wsBF = $wb.Worksheets.Item(6)

Set-TextValue $wsBF "A2" "1.0790051941172372"
Set-TextValue $wsBF "A3" "1.4739279896017703"
